# Apply updated crypto price/volume data per Sat Jul  1 09:36:18 UTC 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.439.03'
$ws.Range('E2').Value = '  -1.22%  '
$ws.Range('D3').Value = '1.917.41'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9994'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4701'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2861'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06835'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '110.37'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.00%  '
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07736'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.71%  '
$ws.Range('D13').Value = '1.890.48'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.304'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6592'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '295.83'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.41%  '
$ws.Range('D17').Value = '30.434.83'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007636'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9996'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('D21').Value = '2.139.87'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9990'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.257'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.226'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '21.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.365'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.089'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1070'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.183'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.993'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05050'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.156'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7369'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02067'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.742'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.677'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.059'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '109.61'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8730'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.834'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4262'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9991'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '51.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +20.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.199'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.234'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1219'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2463'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.84%  '
